$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold, border, centered) from H1 into the
# new I1:J1 header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-30
$data = @{
    2  = @(4, 6)
    3  = @(1, 5)
    4  = @(3, 4)
    5  = @(1, 2)
    6  = @(1, 5)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 5)
    10 = @(1, 7)
    11 = @(1, 2)
    12 = @(1, 5)
    13 = @(1, 5)
    14 = @(1, 4)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 7)
    19 = @(1, 6)
    20 = @(9, 9)
    21 = @(6, 9)
    22 = @(1, 4)
    23 = @(6, 8)
    24 = @(1, 6)
    25 = @(6, 8)
    26 = @(1, 4)
    27 = @(4, 7)
    28 = @(6, 7)
    29 = @(6, 7)
    30 = @(5, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
